$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet from "My Series" to "Data"
$ws.Name = "Data"

# 2) Update cell A11 text
$ws.Range("A11").Value = "Function Information"

# 3) Update cell B20 numeric value (precision correction)
$ws.Range("B20").Value = -0.5631136150295863

# 4) Update the number format used by the data block (B27:C356) from "0.000" to "###0.000"
$ws.Range("B27:C356").NumberFormat = "###0.000"

# 5) Replace the big CEIC metadata blob stored in the A1 cell comment
$newCommentText = "CykAAB+LCAAAAAAAAAPtWltvG8cV/isLAgFaINTuklQkEeMNeJFcoqRkiHRt5yVY7o7IqZe77M6sJL6lL7UQJ0UDw2nd2ikCOGiTxq4Bp41rN/4xgSjJT/kLPTOzV15k0nUfGjgQYs65zZkzZ875dnbRu4cDR9nHPiWeeyGnr2g5BbuWZxO3dyEXsL28/k7uXQNtHlrYuWT65gAzEFZAy6XlQ0ou5PqMDcuqenBwsHJQXPH8nlrQNF292mq2rT4emHniUma6Fs7FWvbLtXIGqtmDFmambTJTal7INdqNlRomVh1oLdM1e9hfqQaUuJjSTZcRRjDlmj42Ga7VW7+QCzMKK++s6EidoieS1YA4tpTLSEp6KAfT4g4ZYKOg6et5bSOvFTu6XtZK5WJxZXVNey9SjAVR06Ssjf19YglCm5mDoVDXNrSirmulYhGpM4XAVhIAA+049i7eJxTbNew4dKmIqOEGViwGq14umBpSU7qhoVd34aJvDvsdwhy8nBtvKUPsKy3PZf23ld1WVRm4oWOJRQNteT62IJiv5N82Ptjxwxh3hk3gdvrEZ6O6OVra1mWK/Z0hj9hyqgaqwxIrDvbZ5SFsPLYhL4BhMD/ASJ3DTJTqhFrwm7gBto0906FppQwTXfH863RoWngbDrXKbRy4jmfakH2MUEasZNIpBrrke7Af3OWq59hbYDUUnsGILTdcCDGftup51xPvZjGR2FWxv7CnA5NF4lN01O57BzuuM2oHXWr5pIvtejWSnslD/HSG2rWAMm8AXiQkJGkpim2rrZY6gv/gTE4yUR1bZGA6lxwIJTXgNGcJqBIwb4+wmucEA5dGnk1Q0RVYVwcfxuuMx2gHttjloffchhvJy2DPZGUVdr2DeM5phghFilyhVrTp04xJ4TrQok2c5oh94avcIg60jPSOpKjZ3Gj3MWYzE0NyEK+OW7wJGdXRdjDowiHrwknbF7NSpCZ8BNkKGQ9+GRo0lrz462haWfyBHzEbbbr2fLmIiWC61FyGvorUCRKCNTlVx3SvA/UKYf3tSrSWGRwkIzBXfpqH4PwOHXMkyHGU0jTUcC0nsLEsCw13T6Qo901u6lw2miI14aQbyHRHndEQSjUlZQY/LuSgeZcp8wEe5AzLC1zmj3j9QGoo+jIdGnRdMYHpLKyz5+NfBYBKRluBa9U8e/HZbBmdyy5hi3voBb4siouriOjx+hjQOuaVRpT+hfWtZdZE/aXEBy4eeC6xFo82BJl7b7/CQmh0qhbWwPJ8LSzvQHeXzY+f9YXVfICU0O6WmqZCqWcRkazh8bBT+uqcI1PHe2bgAJxj0Gh7ce2dJKMKvT4pkyahy74TVUCDg2UKaNmyBysWQAiOCFcsb8AJKoDUK22kpuU5DrLwpttrmm4vAKQR15VJelx/eZfs+KZL+XJiYDFRimcLoahOScBjyOK1E4hEkMXLAy5SJ+RQBw+Gnm86LQgM2QrTLkRNAEhaJuuHI+htDraiIKuJaqyV9Sxy/GVioknJZfADH5bJCaIQ4muRuDyRSWiIr7IFx9KpmQ7p+rKqRq18Fg82LIGIUf3li1sSLkZ7AM9m0H1/jkccsCeDkC5SVo8YMoF5ITXau6X1QmljHbqZGCOx4joeepQwZRcWWFZqfQJovQaKEEmlCntYVir7vbKyRVx4riNAbI8ow4OyUtQU5inrGwrHzNAgRfjCFvVaLadtoq2oP4gl1GEtWUpWAHBKD8w704IxJ1Ew6iZxRik5GaGmZ4HY+MEfTr5+Pj768vS3N3749x9P7jyF38dPHozvP3xx659nn38ExPG/Ho/vce6LG5+cffbx8Xe3Th8/g2FRO7vxeH1jfP9LGSZpEXXMroOF053q+rpWLEEOxyTEN04V4NsOLCZo164JzB2PUfgcKQa1zUbtYrMqalVMjNRlu1L5I+rIC5JhW65UTCTSRY2yTIoYnaj2heMMN9X+DP7IuI+z0mn+PEUZi9Nnfzl99mCudhiwBMfpGxureQ3+CufDPHj4LnE5fWMWzGtm+kssXEgbnZBBu7K5xHFq2EZRh+f71TVNj/uEHR+EWUKTrNBSx+ypE3qSVJPQK06B9DhiisPRwZTFbHlcUgMZxJN/3Dz7++2MVBjdkJK1As4JiMQnU6OBML2921HaO5d3a5tKZ7PN8yThpeSk8XOEw9njQ5dOGnEg31Zyq4qXYO9c6oRmMmwWVVpfxs6kPxd9LxjKYpNSSKgzJJOqM0tjRuURPFkzJ0pQwpkhHe7o77+YpRCuo56g4vS9SpyFgoGm2ZKeEgonu3ska+BbGQvhZPE1A2Q1nJ30ME5yKHBh55qgoCttEdDr2vupThUS+SPZJY+4jBqltVVdPI+FYwTKOrcn/kWNAbRPYVqEDegTFPQzk24esvAgG9tIzRLA06EJndtLnmNjgqzZSXhf3Pvs5E/fnHz66OzGV+MP/zr+6NPTZ38+e/i5PGUntx+dfPwwrOqThV/4wp+OJaBUxI2LpfDTp3AcoHz/wS3F9ZgC8EUJRAX6/oM7KWPcUQF0EssAD2NHsi5MiaaVuZ6SciX2IaMXq0gwUeMtqxhLhE3LGxIrmeS9PDcFGFcRjJ80OvmAYsUDZPZTWElWOFFeVC9UkS300ppW0AshV3rDl9A1aSr0Fx2vC0AjYojLjAmRjNb5ComsmO9ic6daaSYi0okd38Y+T0P5A0XwlLeQBo1GUaqlKMAFEGkFDr9/mhKbZsWWU8VMDS9y9io2v+Wffe2RkUC1wPclSHLDVwXtYAjAOr6rmS8gbkJTYHpbAt80vE7GjXqWD+MUFzpfls0Jgi8KU8iSRapB+aWRxMbbPDbJEHiZ21OIR/g6QEKrfQJ7ovLSs+n7nj+z/iScSKwFsBxKipqEPJYRmyohvJ1sVkSIat7rgfKrWmktA+W3Ko2yMmJ2WWm4NoG1BRxph88RGWLN4ze3poRs08j9FQ2dA9RbC+L0SblZMF10qLlAHXwXUPyb8dE9gN3Hzz6E4fGTrwCNj7/9Ynz7NydH3wLl9ObXp3+7Of7d0fF3d0/vP38DzF8PMC9oWmmR+1eJtUtz7l9fJzBfLy0AzKXQG2C+KDCfPoILIvFzFedD79bS0HumxrnQe9q1RcD33aMFwXf0HnMad0ecOZD79M7z4yc3jp8+PXn0yfjJr/9L3K3Px936+6m2Mo279ULpDex+A7v/r2C3/iOA3fIO+kePuieqz/8EdavJ1XtcGL06djBb7gsVNdFuefuvrAt7v6xqg+44dhjM5d4exGFJDKQ/0+GJ8lq/0pGZV/F9wEz8hf7Sn9VEL6J2Tbe3pGtyXUKRv9GB2cPvXraIT9lV3p3CX5JyLaZck+DzqlGQCPOqHF8ziqurkgISatq8mvEzOshMfjrlOU0yIEu+59Gi0541AsEcDiU8ayyXN7zRbONDAI8pC1Aiu7+EJiJfjC5jTaYvVNZYn3+MQEmvz5Z1bK1rYht3tbzVxYV8ydbW8xsYF/O6Dv83rQI8TqzyTxlC41BHCD5YchI12rDkk0bjP1zgAB4LKQAA"
$comment = $ws.Range("A1").Comment
[void]$comment.Text($newCommentText)
